$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164169430732727
$ws.Range("B1").Value = 2.4220130443573
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.376243591308594
$ws.Range("E1").Value = 1.234835505485535
